{"js": "// Replace the date line and every \"NN\u00d7NN=\" multiplication prompt in the\n// table with the new values from the commit. Every old string is unique\n// within the document, so a scoped exact-text search + replace is safe.\nconst replacements = [\n  [\"2025-09-26 Friday\", \"2025-09-27 Saturday\"],\n  [\"51\u00d734=\", \"40\u00d724=\"],\n  [\"52\u00d753=\", \"77\u00d799=\"],\n  [\"75\u00d731=\", \"79\u00d746=\"],\n  [\"54\u00d778=\", \"55\u00d779=\"],\n  [\"38\u00d795=\", \"17\u00d725=\"],\n  [\"75\u00d759=\", \"81\u00d781=\"],\n  [\"98\u00d762=\", \"70\u00d714=\"],\n  [\"91\u00d744=\", \"56\u00d720=\"],\n  [\"20\u00d766=\", \"13\u00d740=\"],\n  [\"95\u00d732=\", \"13\u00d769=\"],\n  [\"18\u00d723=\", \"74\u00d780=\"],\n  [\"46\u00d755=\", \"67\u00d732=\"],\n  [\"47\u00d794=\", \"90\u00d759=\"],\n  [\"84\u00d725=\", \"41\u00d749=\"],\n  [\"36\u00d750=\", \"74\u00d794=\"],\n  [\"37\u00d715=\", \"81\u00d757=\"],\n  [\"58\u00d780=\", \"28\u00d744=\"],\n  [\"39\u00d713=\", \"94\u00d741=\"],\n  [\"67\u00d738=\", \"22\u00d761=\"],\n  [\"15\u00d793=\", \"97\u00d779=\"],\n  [\"43\u00d755=\", \"98\u00d720=\"],\n  [\"35\u00d713=\", \"53\u00d753=\"],\n  [\"72\u00d721=\", \"34\u00d790=\"],\n  [\"14\u00d762=\", \"37\u00d742=\"],\n  [\"78\u00d750=\", \"67\u00d767=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every \"NN\u00d7NN=\" multiplication prompt in the\n# table with the new values from the commit. Every old string is unique\n# within the document, so Find/Replace (ReplaceAll) for each exact pair\n# is safe and won't touch unrelated text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-09-26 Friday\", \"2025-09-27 Saturday\"),\n    @(\"51\u00d734=\", \"40\u00d724=\"),\n    @(\"52\u00d753=\", \"77\u00d799=\"),\n    @(\"75\u00d731=\", \"79\u00d746=\"),\n    @(\"54\u00d778=\", \"55\u00d779=\"),\n    @(\"38\u00d795=\", \"17\u00d725=\"),\n    @(\"75\u00d759=\", \"81\u00d781=\"),\n    @(\"98\u00d762=\", \"70\u00d714=\"),\n    @(\"91\u00d744=\", \"56\u00d720=\"),\n    @(\"20\u00d766=\", \"13\u00d740=\"),\n    @(\"95\u00d732=\", \"13\u00d769=\"),\n    @(\"18\u00d723=\", \"74\u00d780=\"),\n    @(\"46\u00d755=\", \"67\u00d732=\"),\n    @(\"47\u00d794=\", \"90\u00d759=\"),\n    @(\"84\u00d725=\", \"41\u00d749=\"),\n    @(\"36\u00d750=\", \"74\u00d794=\"),\n    @(\"37\u00d715=\", \"81\u00d757=\"),\n    @(\"58\u00d780=\", \"28\u00d744=\"),\n    @(\"39\u00d713=\", \"94\u00d741=\"),\n    @(\"67\u00d738=\", \"22\u00d761=\"),\n    @(\"15\u00d793=\", \"97\u00d779=\"),\n    @(\"43\u00d755=\", \"98\u00d720=\"),\n    @(\"35\u00d713=\", \"53\u00d753=\"),\n    @(\"72\u00d721=\", \"34\u00d790=\"),\n    @(\"14\u00d762=\", \"37\u00d742=\"),\n    @(\"78\u00d750=\", \"67\u00d767=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replace, [ref]2)\n}\n"}
